$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Country name reordering in shared strings (swap pairs) ---
$ws.Range("A16").Value = "Arabia Saudita"
$ws.Range("A17").Value = "Pakistan"
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("A80").Value = "Estado de Palestina"

# --- Updated timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 15:24"

# --- Updated statistics values ---
$ws.Range("B4").Value = 4864222
$ws.Range("C4").Value = 2048
$ws.Range("E4").Value = 2256916
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 159001
$ws.Range("B6").Value = 1865947
$ws.Range("C6").Value = 10616
$ws.Range("D6").Value = 1237885
$ws.Range("E6").Value = 588981
$ws.Range("G6").Value = 110
$ws.Range("H6").Value = 39081
$ws.Range("B16").Value = 281456
$ws.Range("C16").Value = 1363
$ws.Range("D16").Value = 243713
$ws.Range("E16").Value = 34759
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 2984
$ws.Range("B17").Value = 280461
$ws.Range("C17").Value = 432
$ws.Range("D17").Value = 249397
$ws.Range("E17").Value = 25065
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 5999
$ws.Range("D22").Value = 94129
$ws.Range("E22").Value = 108751
$ws.Range("G22").Value = 50
$ws.Range("H22").Value = 3863
$ws.Range("B24").Value = 134722
$ws.Range("C24").Value = 2836
$ws.Range("D24").Value = 96103
$ws.Range("E24").Value = 33602
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 5017
$ws.Range("B28").Value = 111538
$ws.Range("C28").Value = 216
$ws.Range("D28").Value = 108254
$ws.Range("E28").Value = 3107
$ws.Range("B34").Value = 81181
$ws.Range("C34").Value = 37
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 5747
$ws.Range("B42").Value = 68250
$ws.Range("C42").Value = 84
$ws.Range("D42").Value = 63163
$ws.Range("E42").Value = 4513
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 574
$ws.Range("B44").Value = 55955
$ws.Range("C44").Value = 485
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 6150
$ws.Range("B47").Value = 51681
$ws.Range("C47").Value = 112
$ws.Range("D47").Value = 37318
$ws.Range("E47").Value = 12624
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 1739
$ws.Range("B61").Value = 26804
$ws.Range("C61").Value = 738
$ws.Range("D61").Value = 17838
$ws.Range("E61").Value = 8803
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 163
$ws.Range("B62").Value = 26738
$ws.Range("C62").Value = 287
$ws.Range("E62").Value = 12086
$ws.Range("G62").Value = 7
$ws.Range("H62").Value = 605
$ws.Range("B67").Value = 21481
$ws.Range("C67").Value = 96
$ws.Range("D67").Value = 19336
$ws.Range("E67").Value = 1426
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 719
$ws.Range("B79").Value = 12856
$ws.Range("C79").Value = 394
$ws.Range("D79").Value = 6592
$ws.Range("E79").Value = 5891
$ws.Range("G79").Value = 11
$ws.Range("H79").Value = 373
$ws.Range("B80").Value = 12770
$ws.Range("C80").Value = 229
$ws.Range("D80").Value = 6419
$ws.Range("E80").Value = 6267
$ws.Range("H80").Value = 84
$ws.Range("B84").Value = 11202
$ws.Range("C84").Value = 74
$ws.Range("D84").Value = 7108
$ws.Range("E84").Value = 3589
$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 505
$ws.Range("B100").Value = 5318
$ws.Range("C100").Value = 24
$ws.Range("D100").Value = 4517
$ws.Range("E100").Value = 647
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 154
$ws.Range("B143").Value = 1216
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 698
